$p = $ppt.ActivePresentation
$nm = $p.NotesMaster
$cs = $nm.ColorScheme
for ($i=1; $i -le 8; $i++) {
    $c = $cs.Colors($i)
    Write-Host "Color $i RGB:" $c.RGB
}
Write-Host "--- Slide master ---"
$sm = $p.SlideMaster
$cs2 = $sm.ColorScheme
for ($i=1; $i -le 8; $i++) {
    $c = $cs2.Colors($i)
    Write-Host "Color $i RGB:" $c.RGB
}
